$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3457.4707
$ws.Range("I116").Value = 2706.6365
$ws.Range("K116").Value = 2706.6365
$ws.Range("M116").Value = 735.3634999999999
$ws.Range("H132").Value = 6291485.5
$ws.Range("I132").Value = 8335444
$ws.Range("J132").Value = 2381.6924
$ws.Range("K132").Value = 25006332
$ws.Range("L132").Value = 7145.0772
$ws.Range("M132").Value = -25003802
$ws.Range("N132").Value = -12205.0772
$ws.Range("H135").Value = 217.23529
$ws.Range("I135").Value = 174.36363
$ws.Range("K135").Value = 1569.27267
$ws.Range("M135").Value = 965.7273299999999
$ws.Range("H137").Value = 1761.4412
$ws.Range("I137").Value = 1538.2609
$ws.Range("J137").Value = 2228.0908
$ws.Range("K137").Value = 4614.7827
$ws.Range("L137").Value = 6684.2724
$ws.Range("M137").Value = -2064.7827
$ws.Range("N137").Value = -11784.2724
$ws.Range("H138").Value = 1037185.5
$ws.Range("I138").Value = 2015.8334
$ws.Range("J138").Value = 1209713.8
$ws.Range("K138").Value = 6047.5002
$ws.Range("L138").Value = 3629141.4
$ws.Range("M138").Value = -907.5002000000004
$ws.Range("N138").Value = -3639421.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3539.8914
$ws.Range("I32").Value = 3940.054
$ws.Range("J32").Value = 1894.7778
$ws.Range("K32").Value = 3940.054
$ws.Range("L32").Value = 1894.7778
$ws.Range("M32").Value = -3653.054
$ws.Range("N32").Value = -2468.7778
$ws.Range("H61").Value = 1800.3
$ws.Range("I61").Value = 1427
$ws.Range("K61").Value = 1427
$ws.Range("M61").Value = -1215
$ws.Range("H74").Value = 1077.3334
$ws.Range("I74").Value = 1234.8889
$ws.Range("K74").Value = 1234.8889
$ws.Range("M74").Value = -360.8888999999999
$ws.Range("H77").Value = 1077.3334
$ws.Range("I77").Value = 1234.8889
$ws.Range("K77").Value = 6174.4445
$ws.Range("M77").Value = -1806.4445
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H132").Value = 3194.4707
$ws.Range("I132").Value = 2945.6538
$ws.Range("K132").Value = 8836.9614
$ws.Range("M132").Value = -6306.9614
$ws.Range("H136").Value = 1800.3
$ws.Range("I136").Value = 1427
$ws.Range("K136").Value = 4281
$ws.Range("M136").Value = -1731

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 344.75
$ws.Range("I22").Value = 344.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 344.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -171.75
$ws.Range("N22").ClearContents()
$ws.Range("H105").Value = 142858780
$ws.Range("I105").Value = 166668260
$ws.Range("J105").Value = 1950
$ws.Range("K105").Value = 166668260
$ws.Range("L105").Value = 1950
$ws.Range("M105").Value = -166666513
$ws.Range("N105").Value = -5444
$ws.Range("H134").Value = 5924.048
$ws.Range("I134").Value = 1125.8948
$ws.Range("K134").Value = 3377.6844
$ws.Range("M134").Value = -842.6844000000001
$ws.Range("H135").Value = 30780
$ws.Range("J135").Value = 30780
$ws.Range("L135").Value = 30780
$ws.Range("N135").Value = -40920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 446
$ws.Range("I22").Value = 325
$ws.Range("J22").Value = 526.6667
$ws.Range("K22").Value = 325
$ws.Range("L22").Value = 526.6667
$ws.Range("M22").Value = 25
$ws.Range("N22").Value = -1226.6667
$ws.Range("H31").Value = 1095.8267
$ws.Range("I31").Value = 748.8461
$ws.Range("J31").Value = 1880.3043
$ws.Range("K31").Value = 748.8461
$ws.Range("L31").Value = 1880.3043
$ws.Range("M31").Value = -453.8461
$ws.Range("N31").Value = -2470.3043
$ws.Range("H34").Value = 1095.8267
$ws.Range("I34").Value = 748.8461
$ws.Range("J34").Value = 1880.3043
$ws.Range("K34").Value = 748.8461
$ws.Range("L34").Value = 1880.3043
$ws.Range("M34").Value = -546.8461
$ws.Range("N34").Value = -2284.3043
$ws.Range("H99").Value = 2633446.2
$ws.Range("I99").Value = 3291306
$ws.Range("K99").Value = 3291306
$ws.Range("M99").Value = -3289808
$ws.Range("H126").Value = 2633446.2
$ws.Range("I126").Value = 3291306
$ws.Range("K126").Value = 9873918
$ws.Range("M126").Value = -9871448
$ws.Range("H134").Value = 6803920
$ws.Range("I134").Value = 8773004
$ws.Range("K134").Value = 26319012
$ws.Range("M134").Value = -26316477
$ws.Range("H135").Value = 69000
$ws.Range("J135").Value = 69000
$ws.Range("L135").Value = 69000
$ws.Range("N135").Value = -79140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 677631.3
$ws.Range("I4").Value = 449744.5
$ws.Range("K4").Value = 1349233.5
$ws.Range("M4").Value = -1349121.5
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H86").Value = 1500
$ws.Range("I86").Value = 1500
$ws.Range("K86").Value = 4500
$ws.Range("M86").Value = -3314
$ws.Range("H89").Value = 1500
$ws.Range("I89").Value = 1500
$ws.Range("K89").Value = 13500
$ws.Range("M89").Value = -7572
$ws.Range("H107").Value = 8548.77
$ws.Range("I107").Value = 750.6
$ws.Range("K107").Value = 2251.8
$ws.Range("M107").Value = -331.8000000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3117.2632
$ws.Range("I126").Value = 1949.7142
$ws.Range("K126").Value = 5849.142599999999
$ws.Range("M126").Value = -3379.142599999999
$ws.Range("H132").Value = 2707.0715
$ws.Range("I132").Value = 2455.1
$ws.Range("J132").Value = 3337
$ws.Range("K132").Value = 7365.299999999999
$ws.Range("L132").Value = 10011
$ws.Range("M132").Value = -4835.299999999999
$ws.Range("N132").Value = -15071

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8489.733
$ws.Range("I136").Value = 10551.818
$ws.Range("K136").Value = 31655.454
$ws.Range("M136").Value = -29105.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1121.1666
$ws.Range("I136").Value = 434.75
$ws.Range("J136").Value = 1464.375
$ws.Range("K136").Value = 1304.25
$ws.Range("L136").Value = 4393.125
$ws.Range("M136").Value = 1245.75
$ws.Range("N136").Value = -9493.125
